$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header A1 from "Metric" to "KPI"
$ws.Range("A1").Value = "KPI"

# Update the Lost Time row (row 4) B/C values and bump their number format
# to show one decimal place (0.0%) instead of whole-percent (0%).
$ws.Range("B4").Value = 0.045
$ws.Range("C4").Value = 0.035
$ws.Range("B4:C4").NumberFormat = "0.0%"

# Update the active selection to I9
$ws.Range("I9").Select()
